# Apply the cryptocurrency price/volume updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values look numeric (e.g. "0.488", "22.09") but must
# remain stored as text, matching the original inlineStr cells. Force the
# number format to Text before assigning, then restore the default "Normal"
# style afterwards so no stray cell formatting is introduced.
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D13", "D15", "D16", "D17", "D18", "D22", "D25", "D33", "D34", "D35", "D36", "D39", "D40", "D42", "D43", "D44", "D45", "D47", "D48")
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.911.23'
$ws.Range('D3').Value = '1.549.51'
$ws.Range('D5').Value = '206.24'
$ws.Range('D6').Value = '0.488'
$ws.Range('D8').Value = '22.09'
$ws.Range('D13').Value = '1.551.83'
$ws.Range('D15').Value = '0.518'
$ws.Range('D16').Value = '26.901.80'
$ws.Range('D17').Value = '61.61'
$ws.Range('D18').Value = '217.26'
$ws.Range('D22').Value = '4.06'
$ws.Range('D25').Value = '154.04'
$ws.Range('D33').Value = '1.422.15'
$ws.Range('D34').Value = '3.09'
$ws.Range('D35').Value = '1.59'
$ws.Range('D36').Value = '0.969'
$ws.Range('D39').Value = '0.526'
$ws.Range('D40').Value = '0.809'
$ws.Range('D42').Value = '5.70'
$ws.Range('D43').Value = '2.32'
$ws.Range('D44').Value = '0.998'
$ws.Range('D45').Value = '64.54'
$ws.Range('D47').Value = '1.684.25'
$ws.Range('D48').Value = '87.41'

foreach ($cell in $priceCells) {
    $ws.Range($cell).Style = "Normal"
}

# Column E ("Volume(1h)") percentage-change values stay text naturally
# because of the surrounding spaces, so a plain assignment is sufficient.
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('E4').Value = '  -0.39%  '
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('E8').Value = '  +2.75%  '
$ws.Range('E9').Value = '  -0.45%  '
$ws.Range('E10').Value = '  +0.71%  '
$ws.Range('E11').Value = '  -0.41%  '
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('E13').Value = '  -0.30%  '
$ws.Range('E14').Value = '  +0.62%  '
$ws.Range('E15').Value = '  +0.70%  '
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('E17').Value = '  -0.10%  '
$ws.Range('E18').Value = '  +1.28%  '
$ws.Range('E19').Value = '  +1.29%  '
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('E22').Value = '  +0.39%  '
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('E24').Value = '  -0.68%  '
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('E26').Value = '  -0.62%  '
$ws.Range('E27').Value = '  +0.43%  '
$ws.Range('E28').Value = '  +0.80%  '
$ws.Range('E30').Value = '  +1.82%  '
$ws.Range('E31').Value = '  -0.72%  '
$ws.Range('E32').Value = '  -0.57%  '
$ws.Range('E33').Value = '  +3.95%  '
$ws.Range('E34').Value = '  +4.70%  '
$ws.Range('E35').Value = '  +2.50%  '
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('E38').Value = '  +0.74%  '
$ws.Range('E39').Value = '  +0.52%  '
$ws.Range('E40').Value = '  +0.09%  '
$ws.Range('E41').Value = '  -0.39%  '
$ws.Range('E42').Value = '  +3.28%  '
$ws.Range('E43').Value = '  +3.09%  '
$ws.Range('E44').Value = '  +0.56%  '
$ws.Range('E45').Value = '  +1.41%  '
$ws.Range('E46').Value = '  +0.99%  '
$ws.Range('E47').Value = '  -0.36%  '
$ws.Range('E49').Value = '  +1.22%  '
$ws.Range('E50').Value = '  +5.07%  '
$ws.Range('E51').Value = '  +0.42%  '
